$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-03-09"

# Update the label for the March row to reflect the new "through" date
$ws.Range("A4").Value = "March (through 03-09)"

# Update the March row (row 4) values for columns C..I (B4 / 2015 is unchanged)
$ws.Range("C4").Value = 13
$ws.Range("D4").Value = 14
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 11
$ws.Range("G4").Value = 19
$ws.Range("H4").Value = 27
$ws.Range("I4").Value = 41

# Update the Total row (row 5) values for columns C..I (B5 / 2015 is unchanged)
$ws.Range("C5").Value = 100
$ws.Range("D5").Value = 145
$ws.Range("E5").Value = 153
$ws.Range("F5").Value = 90
$ws.Range("G5").Value = 160
$ws.Range("H5").Value = 369
$ws.Range("I5").Value = 342
